$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (mean)
$ws.Range("B3").Value = 9729.093657620586
$ws.Range("D3").Value = 470.3620931978388

# Row 4 (std)
$ws.Range("B4").Value = 4265.946504806998
$ws.Range("D4").Value = 379.4299106840738

# Row 5 (min)
$ws.Range("B5").Value = 1082.013

# Row 6 (25%)
$ws.Range("B6").Value = 6416.017

# Row 7 (50%)
$ws.Range("B7").Value = 8613.699500000006
$ws.Range("D7").Value = 489.001

# Row 8 (75%)
$ws.Range("B8").Value = 13327.13650000001
$ws.Range("D8").Value = 645

# Row 9 (max)
$ws.Range("B9").Value = 23299.14500000001
$ws.Range("D9").Value = 3720.002

# Row 10 (Total)
$ws.Range("F10").Value = 5113611626.446007

# Row 11 (Residential)
$ws.Range("G11").Value = 0.8247830003711643

# Row 12 (Community)
$ws.Range("F12").Value = 247222316.1850001
$ws.Range("G12").Value = 0.04834593126048981

# Row 13 (IGA)
$ws.Range("G13").Value = 0.1268710683683459
